$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 10:51:59"
$wsZhCn.Range("H2").Value = "2016-03-20 10:52:18"

# de-de sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 10:52:02"
$wsDeDe.Range("H2").Value = "2016-03-20 10:52:23"
